$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (prefix sum): attempts/completed bumped, date advanced to Jan 9 ---
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 45666

# --- Row 10 (dp): one attempt logged today ---
$ws.Range("B10").Value = 1

# --- Clear the other Jan-8 date stamps that were not touched today ---
$ws.Range("E2").Clear()
$ws.Range("E3").Clear()
$ws.Range("E5").Clear()

# --- Insert a new "topo sort" row above "bfs" (was row 15) ---
$ws.Rows.Item(15).Insert()
$ws.Range("A15:E15").Clear()
$ws.Range("A15").Value = "topo sort"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("B15:C15").HorizontalAlignment = -4108

# --- Highlight "bfs" and "dfs" rows (now 16 and 17) in blue, log today's dfs attempt ---
$ws.Range("A16:C16").Interior.Color = 15773696
$ws.Range("A17:C17").Interior.Color = 15773696
$ws.Range("B16:C17").HorizontalAlignment = -4108

$ws.Range("E16").NumberFormat = "[$-1010409]d\ mmmm\ yyyy;@"
$ws.Range("E16").Interior.Color = 15773696
$ws.Range("E16").ClearContents()

$ws.Range("E17").Interior.Color = 15773696
$ws.Range("E17").Value = 45666

# --- Move the active selection, matching the author's last click ---
$ws.Range("C9").Select()
